# Apply the cryptos list update (cell content edits only; no structural changes)
# Commit: "Updated cryptos list on Sun Mar 10 13:21:26 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.763.71"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "3.929.66"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "530.75"
$ws.Range("E5").Value = "  +9.15%  "
$ws.Range("D6").Value = "144.89"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.728"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +4.42%  "
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "42.53"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "4.553.46"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "10.35"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "3.943.26"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "1.23"
$ws.Range("E16").Value = "  +8.00%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.136"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "19.88"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "69.572.95"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "436.67"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "14.45"
$ws.Range("E23").Value = "  -4.90%  "
$ws.Range("D24").Value = "4.13"
$ws.Range("E24").Value = "  +13.11%  "
$ws.Range("D25").Value = "88.19"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("D28").Value = "36.65"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "697.30"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").Value = "13.27"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "70.89"
$ws.Range("E32").Value = "  +16.36%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.84"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").Value = "0.450"
$ws.Range("E34").Value = "  +14.23%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("D36").Value = "40.39"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "0.0₃0848"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "0.0483"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").Value = "3.09"
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").Value = "  -7.45%  "
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").Value = "  +12.75%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0345"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "145.39"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").Value = "2.08"
$ws.Range("E51").Value = "  -2.34%  "
